$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 175, shifting existing rows 175-235 down to 176-236.
$ws.Rows.Item(175).Insert()

# Populate the new row 175 with the same values as the (old) row 175 template,
# but with updated Fecha/Volumen/Precio values for the new weekly entry.
$ws.Cells.Item(175, 1).Value = 3
$ws.Cells.Item(175, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(175, 3).Value = "Coquimbo"
$ws.Cells.Item(175, 4).Value = 44559
$ws.Cells.Item(175, 4).NumberFormat = $ws.Cells.Item(176, 4).NumberFormat
$ws.Cells.Item(175, 5).Value = 5
$ws.Cells.Item(175, 6).Value = 100112039
$ws.Cells.Item(175, 7).Value = "Ciboulette"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 172
$ws.Cells.Item(175, 11).Value = 1500
$ws.Cells.Item(175, 12).Value = 2000
$ws.Cells.Item(175, 13).Value = 1747
$ws.Cells.Item(175, 14).Value = "`$/docena de atados"
$ws.Cells.Item(175, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(175, 16).Value = 582
$ws.Cells.Item(175, 17).Value = 3
$ws.Cells.Item(175, 18).Value = "Hortaliza"
